# Insert a new data row at row 292 (pushing existing rows 292-389 down to 293-390)
# and populate it with a new weekly price observation, as described by the commit
# "Fruta / hortaliza, semanal".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 292; Excel shifts formatting/content of
# rows 292..389 down to 293..390, and the new blank row 292 inherits the
# formatting (e.g. the date-formatted column D) from the row above it.
$ws.Rows.Item(292).Insert()

# Populate the newly inserted row 292 with the new observation.
$ws.Range("A292").Value = 7
$ws.Range("B292").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C292").Value = "Ñuble"
$ws.Range("D292").Value = 44988
$ws.Range("E292").Value = 16
$ws.Range("F292").Value = 100114013
$ws.Range("G292").Value = "Zanahoria"
$ws.Range("H292").Value = "Sin especificar"
$ws.Range("I292").Value = "Segunda"
$ws.Range("J292").Value = 60
$ws.Range("K292").Value = 5500
$ws.Range("L292").Value = 5500
$ws.Range("M292").Value = 5500
$ws.Range("N292").Value = "`$/saco 20 kilos"
$ws.Range("O292").Value = "Región de Ñuble"
$ws.Range("P292").Value = 275
$ws.Range("Q292").Value = 20
$ws.Range("R292").Value = "Hortaliza"
